$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Nachname) and B (Vorname), rows 2-8
$data = @(
    @("Matumona", "Noe"),
    @("Zillig", "Nicolas"),
    @("Sarman", "Dominik"),
    @("Kohler", "Alina"),
    @("Matumona", "Nina"),
    @("asdf", "Marlene"),
    @("Kohler", "Nina")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Columns C (Note Exakt) and D (Note Gerundet) are no longer populated for rows 2-8
$ws.Range("C2:D8").ClearContents()
